# Updated cryptos list on Sat Jan  6 11:30:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'43.999.14"
$ws.Range("E2").Value = "  -0.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.237.13"
$ws.Range("E3").Value = "  -0.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'306.01"
$ws.Range("E5").Value = "  -4.08%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'94.85"
$ws.Range("E6").Value = "  -5.83%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.23%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -4.69%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'34.81"
$ws.Range("E10").Value = "  -5.61%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -2.89%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  -4.29%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.42%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'2.578.64"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "'2.236.93"
$ws.Range("E15").Value = "  -0.60%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.827"
$ws.Range("E16").Value = "  -3.09%  "

# Row 17 - Chainlink
$ws.Range("E17").Value = "  -4.79%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'43.883.44"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.0₃0960"
$ws.Range("E19").Value = "  -1.94%  "

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.12"
$ws.Range("E20").Value = "  -9.81%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'6.26"
$ws.Range("E21").Value = "  -2.97%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'64.99"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'236.66"
$ws.Range("E23").Value = "  +1.16%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "'2.94"
$ws.Range("E24").Value = "  -5.32%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -5.41%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.07%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -6.10%  "

# Row 28 - InjectiveProtocol
$ws.Range("D28").Value = "'37.70"
$ws.Range("E28").Value = "  -2.96%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  -1.18%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "'5.95"
$ws.Range("E30").Value = "  -2.89%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'19.88"
$ws.Range("E31").Value = "  -1.19%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'153.40"
$ws.Range("E32").Value = "  -4.07%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0801"
$ws.Range("E33").Value = "  -5.33%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'3.24"
$ws.Range("E34").Value = "  +5.93%  "

# Row 35 - WEMIXToken
$ws.Range("E35").Value = "  -3.97%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -0.23%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -6.24%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  -8.35%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "'15.22"
$ws.Range("E39").Value = "  -7.43%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "'3.83"
$ws.Range("E40").Value = "  -8.11%  "

# Row 41 - NEARProtocol
$ws.Range("D41").Value = "'3.35"
$ws.Range("E41").Value = "  -9.22%  "

# Row 42 - VeChain
$ws.Range("D42").Value = "'0.0301"
$ws.Range("E42").Value = "  -4.06%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.23%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'1.724.48"
$ws.Range("E44").Value = "  -2.31%  "

# Row 45 - BitcoinSV
$ws.Range("D45").Value = "'85.50"
$ws.Range("E45").Value = "  +5.53%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  -4.01%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'100.07"
$ws.Range("E47").Value = "  -3.45%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "'4.91"
$ws.Range("E48").Value = "  -4.89%  "

# Row 49 - was FraxShare, now ordi (rows 49/50 swap identity)
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'68.95"
$ws.Range("E49").Value = "  -7.60%  "

# Row 50 - was ordi, now FraxShare
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.08"
$ws.Range("E50").Value = "  -2.68%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "'54.21"
$ws.Range("E51").Value = "  -5.28%  "
